# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" suffixed column headers (row 1) to
# "_FV2404" / "_FV2410" respectively, wraps the used range in an Excel
# Table (ListObject) with an AutoFilter, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) suffixes ----------------------------
# Columns A:J used the "_old" suffix (FV2404 = the "old" format version),
# columns L:U used the "_new" suffix (FV2410 = the "new" format version).
# Column K ("diff") is left untouched.

$oldHeaders = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$newHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# --- 2. Turn the used range into an Excel Table (ListObject) -------------

$tableRange = $ws.Range("A1:U64")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
# No named table style in the source workbook (plain banded table only).
$lo.TableStyle = ""

# --- 3. Freeze the header row ---------------------------------------------

$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
Write-Output "edit.ps1 applied: headers renamed, table added, header row frozen"
